$wb = $excel.ActiveWorkbook

# --- Rename the first sheet ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Export as TSV"

# --- Freeze the header row (row 1) on that sheet ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Add error alert text to the existing data validations ---
$ws.Range("I2:I1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("I2:I1048576").Validation.ErrorMessage = "Value must be one of: imaging."

$ws.Range("J2:J1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("J2:J1048576").Validation.ErrorMessage = "Value must be one of: AF."

$ws.Range("L2:L1048576").Validation.ErrorTitle = "Not a boolean"
$ws.Range("L2:L1048576").Validation.ErrorMessage = "The values in this column must be ""TRUE"" or ""FALSE""."

$ws.Range("O2:O1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("O2:O1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("P2:P1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("P2:P1048576").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("Q2:Q1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("Q2:Q1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("R2:R1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("R2:R1048576").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("S2:S1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("S2:S1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("T2:T1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("T2:T1048576").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("U2:U1048576").Validation.ErrorTitle = "Not an integer"
$ws.Range("U2:U1048576").Validation.ErrorMessage = "The values in this column must be integers."
